$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2022" header in G1 — match the existing header formatting (bold+centered,
# same style as C1:F1) and keep it a text label like the other year headers.
# A plain Value assignment of "2022" would auto-coerce to a number, so write
# it as a formula that evaluates to the text "2022" and then collapse that
# formula down to a static value via copy / paste-values.
$ws.Range("G1").Formula = '="2022"'
$ws.Range("G1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4163) | Out-Null

# Rename Czech Republic -> Czechia
$ws.Range("B6").Value = "Czechia"

# Updated 2021 (column F) figures
$ws.Range("F7").Value = 92873
$ws.Range("F8").Value = 51435
$ws.Range("F22").Value = 104847

# New 2022 (column G) figures
$ws.Range("G2").Value = 86005
$ws.Range("G3").Value = 16055
$ws.Range("G4").Value = 13604
$ws.Range("G5").Value = 17889
$ws.Range("G6").Value = 16684
$ws.Range("G7").Value = 99443
$ws.Range("G8").Value = 43404
$ws.Range("G9").Value = 12720
$ws.Range("G10").Value = 9908
$ws.Range("G11").Value = 33892
$ws.Range("G12").Value = 17106
$ws.Range("G13").Value = 27567
$ws.Range("G14").Value = 14816
$ws.Range("G15").Value = 96040
$ws.Range("G16").Value = 9777
$ws.Range("G17").Value = 183119
$ws.Range("G18").Value = 26112
$ws.Range("G19").Value = 8510
$ws.Range("G20").Value = 4257
$ws.Range("G21").Value = 4980
$ws.Range("G22").Value = 119757
$ws.Range("G23").Value = 110410
$ws.Range("G24").Value = 14436
$ws.Range("G25").Value = 82240
$ws.Range("G26").Value = 68666
$ws.Range("G27").Value = 30721
$ws.Range("G28").Value = 34441
